$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 2 (ECs) ---
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 1.298672
$ws.Range("H2").Value = 3.896016
$ws.Range("M2").Value = 2.885873333333334
$ws.Range("N2").Value = 8.657620000000001
$ws.Range("O2").Value = 0.3070415651026022
$ws.Range("P2").Value = 0.3070415651026022
$ws.Range("Q2").Value = 3.747802893546667
$ws.Range("R2").Value = 33.73022604192001
$ws.Range("S2").Value = 0.3070415651026022
$ws.Range("T2").Value = 0.3070415651026022

# --- Update row 3 (FAPs) ---
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 1.298672
$ws.Range("H3").Value = 3.896016
$ws.Range("O3").Value = 0.3368406220840099
$ws.Range("P3").Value = 0.3368406220840099
$ws.Range("Q3").Value = 4.111535380197334
$ws.Range("R3").Value = 37.003818421776
$ws.Range("S3").Value = 0.3368406220840099
$ws.Range("T3").Value = 0.3368406220840099

# --- Update row 4 (MuSCs) ---
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 1.298672
$ws.Range("H4").Value = 3.896016
$ws.Range("M4").Value = 3.327024333333334
$ws.Range("N4").Value = 9.981073
$ws.Range("O4").Value = 0.3539776838580724
$ws.Range("P4").Value = 0.3539776838580724
$ws.Range("Q4").Value = 4.320713345018667
$ws.Range("R4").Value = 38.886420105168
$ws.Range("S4").Value = 0.3539776838580724
$ws.Range("T4").Value = 0.3539776838580724

# --- Add new row 5 (Resolving-Mac) ---
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Tgfa"
$ws.Range("C5").Value = "Erbb2"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 1.298672
$ws.Range("H5").Value = 3.896016
$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.020115
$ws.Range("N5").Value = 0.060345
$ws.Range("O5").Value = 0.002140128955315263
$ws.Range("P5").Value = 0.002140128955315263
$ws.Range("Q5").Value = 0.02612278728
$ws.Range("R5").Value = 0.23510508552
$ws.Range("S5").Value = 0.002140128955315263
$ws.Range("T5").Value = 0.002140128955315263
